$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Flight_Search")

# Update the flight dates (the old dates had already passed and were
# triggering notifications) - bump them forward.
$ws.Range("E2").Value = "Wed, 27 Feb, 2019"
$ws.Range("F2").Value = "Thu, 28 Feb, 2019"
$ws.Range("E3").Value = "Wed, 27 Feb, 2019"

# Move the active selection to F2 to match the saved view state.
$ws.Range("F2").Select()
